# Car Profile Sync.pptx edit script
# - Slide 2 ("Idee"): fix the rotated "0-Click" textbox -> "1-Click", no rotation;
#   add click-triggered entrance (fade) animations for the two car/phone groups.
# - Slide 3 ("Vorteile"): rename "Kontaktlos" -> "Einfachheit"; add click-triggered
#   entrance (fade) animations for the "Sicher" and "Erweiterbar" groups.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Idee")
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Shape 14 = "Textfeld 22" (id="23"), the rotated "   0-Click" textbox nearest
# the second car image. Drop the rotation and split the run so "0" becomes "1".
$lblShape = $s2.Shapes.Item(14)
$lblShape.Rotation = 0

$tr = $lblShape.TextFrame.TextRange
$numPart = $tr.Characters(4, 7)
$numPart.Text = "1-Click"

# Animations: two click-triggered build groups, each an entrance Fade (id 10)
# applied to a "lead" shape on click, with the remaining shapes in that group
# fading in at the same time ("With Previous").

# Group 1: the left car/arrow group, together with the big picture behind it.
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(6), 10)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(7), 10, 0, 2)

# Group 2: the second "1-Click" label group, together with its picture, the
# other arrow group, its label, the last picture and the last arrow group.
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(10), 10)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(4), 10, 0, 2)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(9), 10, 0, 2)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(12), 10, 0, 2)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(8), 10, 0, 2)
$e = $s2.TimeLine.MainSequence.AddEffect($s2.Shapes.Item(11), 10, 0, 2)

# ---------------------------------------------------------------------------
# Slide 3 ("Vorteile")
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Shape 6 = "Textfeld 8" (id="9"), currently reads "Kontaktlos".
$advShape = $s3.Shapes.Item(6)
$advShape.TextFrame.TextRange.Text = "Einfachheit"

# Animations: two click-triggered build groups (entrance Fade, id 10).

# Group 1: the lock picture, together with the "Sicher" caption.
$e = $s3.TimeLine.MainSequence.AddEffect($s3.Shapes.Item(3), 10)
$e = $s3.TimeLine.MainSequence.AddEffect($s3.Shapes.Item(4), 10, 0, 2)

# Group 2: the extra picture, together with the gear picture and the
# "Erweiterbar" caption.
$e = $s3.TimeLine.MainSequence.AddEffect($s3.Shapes.Item(7), 10)
$e = $s3.TimeLine.MainSequence.AddEffect($s3.Shapes.Item(1), 10, 0, 2)
$e = $s3.TimeLine.MainSequence.AddEffect($s3.Shapes.Item(5), 10, 0, 2)
